$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lito-info")

# Insert a new row above row 2, shifting all existing rows down,
# then populate it with the new 2020-21 figures (post-Budget 2019 update).
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "2020-21"
$ws.Range("B2").Value = 445
$ws.Range("C2").Value = 0.015
$ws.Range("D2").Value = 37000
